$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values (row 2) ---
# C2 holds the SQL server IP; fix it from the (broken-on-Linux) loopback
# address to the real server IP.
$ws.Range("C2").Value = "192.168.1.113"
# Match C2's number format to the neighbouring text-formatted cells (B2/A2)
# so it is stored as text, not accidentally parsed as a number.
$ws.Range("C2").NumberFormat = "@"

# --- Column widths: split the merged B:C column width definition so C can
# get its own, wider, best-fit width to accommodate the longer IP text. ---
$ws.Columns.Item(3).ColumnWidth = 14.29

# --- Update the active selection shown when the sheet is opened. ---
$ws.Range("C2").Select()
